# Objective function adapted + minor changes
# - Insert a new "Set" row for "Capacity defining commodity group" (\Cdcg)
#   right above the existing "Output commodities of a specific unit" row,
#   i.e. at worksheet row 13 (pushing all following rows down by one).
# - Simplify the "UnitCapacity" parameter's "optional" column (E) from
#   "\unit,Cdcg" to "\unit" now that Cdcg is its own Set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13; existing row 13 (and everything below it)
# shifts down to row 14, etc.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new Set definition.
$ws.Cells.Item(13, 1).Value = "Set"
$ws.Cells.Item(13, 3).Value = "Capacity defining commodity group"
$ws.Cells.Item(13, 4).Value = "\Cdcg"
$ws.Cells.Item(13, 5).Value = "\commodity,\unit"
$ws.Cells.Item(13, 6).Value = "Cdcg_{#1}"
$ws.Cells.Item(13, 7).Value = "Set of commodities of which the sum of the flows are restricted by the capacity of the unit"
$ws.Cells.Item(13, 9).Value = "s_topology"
$ws.Cells.Item(13, 10).Value = 1
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 1

# The "UnitCapacity" parameter row (now shifted to row 25) no longer
# references the inline "Cdcg" note in its optional-index column;
# it is simplified to just "\unit" now that Cdcg is its own Set.
$ws.Cells.Item(25, 5).Value = "\unit"

# The worksheet table ("Table1") covers the data range; grow it to
# include the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:Q260"))

# Cosmetic: update the active selection to match.
$ws.Range("I13").Select()
